# Updated cryptos list values per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "42.172.98"
$ws.Cells.Item(2, 5).Value = "  +1.41%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.212.54"
$ws.Cells.Item(3, 5).Value = "  -0.37%  "

$ws.Cells.Item(4, 5).Value = "  -0.06%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "231.26"
$ws.Cells.Item(5, 5).Value = "  +1.41%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.613"
$ws.Cells.Item(6, 5).Value = "  -1.10%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "60.89"
$ws.Cells.Item(7, 5).Value = "  +0.07%  "

$ws.Cells.Item(8, 5).Value = "  -0.04%  "

$ws.Cells.Item(9, 5).Value = "  +0.47%  "

$ws.Cells.Item(10, 5).Value = "  +2.60%  "

$ws.Cells.Item(11, 5).Value = "  +0.43%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "2.541.63"
$ws.Cells.Item(12, 5).Value = "  -0.49%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "15.45"
$ws.Cells.Item(13, 5).Value = "  -1.05%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "22.07"
$ws.Cells.Item(14, 5).Value = "  +3.01%  "

$ws.Cells.Item(15, 5).Value = "  +0.38%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "5.57"
$ws.Cells.Item(16, 5).Value = "  +0.63%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "2.224.44"
$ws.Cells.Item(17, 5).Value = "  +0.15%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "42.058.62"
$ws.Cells.Item(18, 5).Value = "  +1.39%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.0₃0935"
$ws.Cells.Item(19, 5).Value = "  +5.97%  "

$ws.Cells.Item(20, 5).Value = "  +3.00%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "71.95"
$ws.Cells.Item(21, 5).Value = "  -0.60%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "243.25"
$ws.Cells.Item(22, 5).Value = "  -1.18%  "

$ws.Cells.Item(23, 5).Value = "  -0.13%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.45"
$ws.Cells.Item(24, 5).Value = "  +4.03%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.36"
$ws.Cells.Item(25, 5).Value = "  -0.34%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "9.58"
$ws.Cells.Item(26, 5).Value = "  +0.53%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "168.90"
$ws.Cells.Item(27, 5).Value = "  +0.58%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.140"
$ws.Cells.Item(28, 5).Value = "  +0.89%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "20.28"
$ws.Cells.Item(29, 5).Value = "  +2.22%  "

$ws.Cells.Item(30, 5).Value = "  +3.13%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "2.68"
$ws.Cells.Item(31, 5).Value = "  +1.37%  "

$ws.Cells.Item(32, 5).Value = "  -0.96%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.96"
$ws.Cells.Item(33, 5).Value = "  -1.78%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "4.60"
$ws.Cells.Item(34, 5).Value = "  -0.59%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.0647"
$ws.Cells.Item(35, 5).Value = "  +4.45%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "6.30"
$ws.Cells.Item(36, 5).Value = "  -4.28%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.55"
$ws.Cells.Item(37, 5).Value = "  -3.64%  "

$ws.Cells.Item(38, 5).Value = "  -1.44%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0248"
$ws.Cells.Item(39, 5).Value = "  +6.14%  "

$ws.Cells.Item(40, 5).Value = "  +0.24%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.000227"
$ws.Cells.Item(41, 5).Value = "  -2.44%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "8.52"
$ws.Cells.Item(42, 5).Value = "  -3.20%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.0955"
$ws.Cells.Item(43, 5).Value = "  -2.13%  "

$ws.Cells.Item(44, 5).Value = "  +1.95%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "96.82"
$ws.Cells.Item(45, 5).Value = "  -1.92%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.457.12"
$ws.Cells.Item(46, 5).Value = "  -0.60%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "4.31"
$ws.Cells.Item(47, 5).Value = "  -10.89%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "16.08"
$ws.Cells.Item(48, 5).Value = "  -1.17%  "

$ws.Cells.Item(49, 2).Value = "ARBITRUM"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.06"
$ws.Cells.Item(49, 5).Value = "  -0.69%  "

$ws.Cells.Item(50, 2).Value = "HuobiToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.71"
$ws.Cells.Item(50, 5).Value = "  -2.64%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.19"
$ws.Cells.Item(51, 5).Value = "  +2.02%  "
